# Update "想去人数" (interest count) values in column F for the matching
# rows on both the "展览" and "全部类型" sheets. Each target row is located
# by its unique event name in column C so the script is robust even if row
# numbers shift slightly between sheets.

$wb = $excel.ActiveWorkbook

# Map of event name -> new F-column value.
$updates = @{
    "北京·井上直久の依巴拉度世界"             = 20
    "北京·EXA·全职高手ONLY·夏令营"            = 520
    "北京·Aw动漫游戏嘉年华8th"                 = 5069
    "北京·IDO动漫游戏嘉年华46th"               = 3730
    "北京·配音演员 刘北辰 专场活动"            = 206
    "北京·第17届IJOY漫展xCGF游戏节"            = 6629
    "北京·IDO暑假狂欢节"                       = 1065
    "北京·万游引力国潮动漫嘉年华s8"            = 500
    "北京·梦次元动漫展M30"                     = 1351
    "北京·第18届IJOY漫展xCGF游戏节"            = 2266
    "北京·IDO动漫游戏嘉年华47th"               = 772
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $name = $ws.Cells.Item($r, 3).Text
        if ($null -ne $name -and $updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value2 = $updates[$name]
        }
    }
}
